$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, E, F, H, J, K for rows 2-13
# Format: row = @(B, C, E, F, H, J, K)
$data = @{
    2  = @(5, 4, 5, 3, 1, 9, 374)
    3  = @(5, 5, 5, 5, 2, 2, 12)
    4  = @(2, 1, 2, 2, 1, 9, 374)
    5  = @(1, 1, 2, 2, 2, 2, 12)
    6  = @(2, 2, 4, 4, 2, 3, 22)
    7  = @(2, 2, 2, 2, 1, 9, 374)
    8  = @(4, 4, 5, 5, 5, 3, 22)
    9  = @(1, 1, 2, 2, 2, 3, 22)
    10 = @(1, 1, 1, 1, 1, 3, 22)
    11 = @(1, 1, 2, 2, 2, 3, 22)
    12 = @(15, 15, 7, 7, 3, 3, 22)
    13 = @(3, 3, 8, 8, 7, 10, 380)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("E$row").Value = $vals[2]
    $ws.Range("F$row").Value = $vals[3]
    $ws.Range("H$row").Value = $vals[4]
    $ws.Range("J$row").Value = $vals[5]
    $ws.Range("K$row").Value = $vals[6]
}
